$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.949.02'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.768.92'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.00'
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4582'
$ws.Range("E7").Value = '  +2.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3526'
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.12'
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07397'
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.096'
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.73'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.008'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.193'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.765.07'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.66'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001061'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06448'
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.773'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.968.41'
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.147'
$ws.Range("E25").Value = '  +2.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.30'
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.18'
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.969.72'
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.163'
$ws.Range("E29").Value = '  +3.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.11'
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.077'
$ws.Range("E31").Value = '  -1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09300'
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.600'
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.666'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.85'
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02280'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06126'
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2087'
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.939'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6262'
$ws.Range("E40").Value = '  -0.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.183'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.379'
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.815'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.13'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.735'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5856'
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.43'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.937'
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.130'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06817'
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.31'
$ws.Range("E51").Value = '  +2.31%  '
